# Update the "Лист1" sheet of the form workbook with a new data row
# (login/password/name/surname/middlename/class), widen the columns,
# move the selection, and set up the page for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New second row of data (values pulled into the shared-string table).
$ws.Range("A2").Value = "ЛогинПользователя"
$ws.Range("B2").Value = "ПарольПользователя"
$ws.Range("C2").Value = "Максим"
$ws.Range("D2").Value = "Масимов"
$ws.Range("E2").Value = "Максимович"
$ws.Range("F2").Value = "9-5"

# Widen columns A:F.
$ws.Range("A1:F2").Columns.ColumnWidth = 26.28515625

# Move the active selection.
$ws.Range("F12").Select()

# Configure the page for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the workbook window to a maximized-looking layout.
$excel.ActiveWindow.WindowState = -4137
